$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize the "Unknown" placeholder text to lowercase "unknown"
# for the importance/experience columns D2:J2 (state.State in K2 is untouched).
$ws.Range("D2:J2").Value = "unknown"
